$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new TextContent rows for logout / user-details feature
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Profile"
$ws.Range("C26").Value = 1

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Log Out"
$ws.Range("C27").Value = 1

# Update the active selection as recorded in the saved workbook
$ws.Range("L11").Select()
